$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1454.5714
$ws.Range("I19").Value = 1340
$ws.Range("K19").Value = 1340
$ws.Range("M19").Value = -1165

$ws.Range("H32").Value = 2461.818
$ws.Range("J32").Value = 3297.4
$ws.Range("L32").Value = 3297.4
$ws.Range("N32").Value = -3949.4

$ws.Range("H40").Value = 1295.8667
$ws.Range("I40").Value = 1279.7693
$ws.Range("J40").Value = 1400.5
$ws.Range("K40").Value = 1279.7693
$ws.Range("L40").Value = 1400.5
$ws.Range("M40").Value = -1104.7693
$ws.Range("N40").Value = -1750.5

$ws.Range("H88").Value = 1515.6154
$ws.Range("J88").Value = 1030.5
$ws.Range("L88").Value = 1030.5
$ws.Range("N88").Value = -1842.5

$ws.Range("H91").Value = 1515.6154
$ws.Range("J91").Value = 1030.5
$ws.Range("L91").Value = 1030.5
$ws.Range("N91").Value = -3838.5

$ws.Range("H98").Value = 778.4091
$ws.Range("I98").Value = 447.42105
$ws.Range("K98").Value = 447.42105
$ws.Range("M98").Value = 1050.57895

$ws.Range("H100").Value = 4668.3335
$ws.Range("I100").Value = 4668.3335
$ws.Range("K100").Value = 4668.3335
$ws.Range("M100").Value = -4127.3335

$ws.Range("H105").Value = 670.25
$ws.Range("J105").Value = 670.25
$ws.Range("L105").Value = 670.25
$ws.Range("N105").Value = -7658.25

$ws.Range("H115").Value = 316.42856
$ws.Range("I115").Value = 316.42856
$ws.Range("K115").Value = 949.28568
$ws.Range("M115").Value = 617.71432

$ws.Range("H122").Value = 778.4091
$ws.Range("I122").Value = 447.42105
$ws.Range("K122").Value = 1342.26315
$ws.Range("M122").Value = 1107.73685

$ws.Range("H138").Value = 1442.1578
$ws.Range("I138").Value = 908.7857
$ws.Range("J138").Value = 2935.6
$ws.Range("K138").Value = 2726.3571
$ws.Range("L138").Value = 8806.799999999999
$ws.Range("M138").Value = 2413.6429
$ws.Range("N138").Value = -19086.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 8000
$ws.Range("I26").Value = 8000
$ws.Range("J26").Value = 8000
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = -7670
$ws.Range("N26").Value = -8660

$ws.Range("H38").Value = 9009.5
$ws.Range("I38").Value = 9009.5
$ws.Range("K38").Value = 9009.5
$ws.Range("M38").Value = -8542.5

$ws.Range("H61").Value = 2346.0605
$ws.Range("I61").Value = 2147.1785
$ws.Range("K61").Value = 2147.1785
$ws.Range("M61").Value = -1935.1785

$ws.Range("H136").Value = 2346.0605
$ws.Range("I136").Value = 2147.1785
$ws.Range("K136").Value = 6441.5355
$ws.Range("M136").Value = -3891.5355

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1201.75
$ws.Range("I80").Value = 118
$ws.Range("J80").Value = 1852
$ws.Range("K80").Value = 118
$ws.Range("L80").Value = 1852
$ws.Range("M80").Value = 880
$ws.Range("N80").Value = -3848

$ws.Range("H83").Value = 1201.75
$ws.Range("I83").Value = 118
$ws.Range("J83").Value = 1852
$ws.Range("K83").Value = 590
$ws.Range("L83").Value = 9260
$ws.Range("M83").Value = 4402
$ws.Range("N83").Value = -19244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2196.9092
$ws.Range("I134").Value = 2351.8333
$ws.Range("J134").Value = 1499.75
$ws.Range("K134").Value = 7055.499899999999
$ws.Range("L134").Value = 4499.25
$ws.Range("M134").Value = -4520.499899999999
$ws.Range("N134").Value = -9569.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 120.75
$ws.Range("I2").Value = 110.5
$ws.Range("K2").Value = 663
$ws.Range("M2").Value = -550

$ws.Range("H26").Value = 1716.1852
$ws.Range("I26").Value = 1904.5
$ws.Range("K26").Value = 5713.5
$ws.Range("M26").Value = -5425.5

$ws.Range("H38").Value = 193.14285
$ws.Range("J38").Value = 64
$ws.Range("L38").Value = 192
$ws.Range("N38").Value = -886

$ws.Range("H120").Value = 10000
$ws.Range("J120").Value = 10000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676

$ws.Range("H131").Value = 1874.75
$ws.Range("I131").Value = 1749.75
$ws.Range("J131").Value = 1999.75
$ws.Range("K131").Value = 5249.25
$ws.Range("L131").Value = 5999.25
$ws.Range("M131").Value = -209.25
$ws.Range("N131").Value = -16079.25

$ws.Range("H132").Value = 1866.3636
$ws.Range("I132").Value = 1825.5555
$ws.Range("K132").Value = 16429.9995
$ws.Range("M132").Value = -13899.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3329.3333
$ws.Range("I126").Value = 2494
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 7482
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -5012
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 3296.72
$ws.Range("I132").Value = 3100.25
$ws.Range("K132").Value = 9300.75
$ws.Range("M132").Value = -6770.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 410.66666
$ws.Range("I16").Value = 337
$ws.Range("K16").Value = 337
$ws.Range("M16").Value = -167

$ws.Range("H22").Value = 800.0454999999999
$ws.Range("I22").Value = 583.4167
$ws.Range("J22").Value = 1060
$ws.Range("K22").Value = 583.4167
$ws.Range("L22").Value = 1060
$ws.Range("M22").Value = -288.4167
$ws.Range("N22").Value = -1650

$ws.Range("H27").Value = 800.0454999999999
$ws.Range("I27").Value = 583.4167
$ws.Range("J27").Value = 1060
$ws.Range("K27").Value = 583.4167
$ws.Range("L27").Value = 1060
$ws.Range("M27").Value = -476.4167
$ws.Range("N27").Value = -1274

$ws.Range("H61").Value = 2168
$ws.Range("I61").Value = 2168
$ws.Range("K61").Value = 2168
$ws.Range("M61").Value = -1966

$ws.Range("H95").Value = 30316
$ws.Range("J95").Value = 30316
$ws.Range("L95").Value = 30316
$ws.Range("N95").Value = -35808

$ws.Range("H99").Value = 16999.666
$ws.Range("I99").Value = 18375
$ws.Range("J99").Value = 14249
$ws.Range("K99").Value = 18375
$ws.Range("L99").Value = 14249
$ws.Range("M99").Value = -15380
$ws.Range("N99").Value = -20239

$ws.Range("H106").Value = 24453.334
$ws.Range("J106").Value = 24453.334
$ws.Range("L106").Value = 24453.334
$ws.Range("N106").Value = -26977.334

$ws.Range("H113").Value = 2168
$ws.Range("I113").Value = 2168
$ws.Range("K113").Value = 2168
$ws.Range("M113").Value = 2

$ws.Range("H132").Value = 1631.75
$ws.Range("I132").Value = 1680.1111
$ws.Range("J132").Value = 1196.5
$ws.Range("K132").Value = 5040.3333
$ws.Range("L132").Value = 3589.5
$ws.Range("M132").Value = -2510.3333
$ws.Range("N132").Value = -8649.5

$ws.Range("H136").Value = 25002734
$ws.Range("I136").Value = 2937.9167
$ws.Range("J136").Value = 62502428
$ws.Range("K136").Value = 8813.750100000001
$ws.Range("L136").Value = 187507284
$ws.Range("M136").Value = -6263.750100000001
$ws.Range("N136").Value = -187512384

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 28946.334
$ws.Range("J46").Value = 28946.334
$ws.Range("L46").Value = 28946.334
$ws.Range("N46").Value = -29408.334

$ws.Range("H92").Value = 37033
$ws.Range("J92").Value = 37033
$ws.Range("L92").Value = 37033
$ws.Range("N92").Value = -42025

$ws.Range("H122").Value = 2936.7273
$ws.Range("I122").Value = 2955.4
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 8866.200000000001
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -6416.200000000001
$ws.Range("N122").Value = -13150

$ws.Range("H134").Value = 28946.334
$ws.Range("J134").Value = 28946.334
$ws.Range("L134").Value = 86839.00199999999
$ws.Range("N134").Value = -91909.00199999999
